$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the old "Hoja" column (was S) ---
# The old S column ("Hoja" header / "POAI_2025" data) is pushed one
# column to the right (T); the new column takes the old S position and
# gets header "Col19" with data value 0 in every row.

# Copy the header formatting (bold + border) from S1 into the new T1
# header cell before writing values, so the inserted column matches the
# existing header look.
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

# Shift the old "Hoja" column data (currently in S) over to T.
$ws.Range("T1").Value = $ws.Range("S1").Value()
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 20).Value = $ws.Cells.Item($r, 19).Value()
}

# New column header + data in S.
$ws.Range("S1").Value = "Col19"
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 19).Value = 0
}

# --- Responsable (L) corrections ---
$ws.Range("L8").Value = "SARA DIANA URBANO"
$ws.Range("L11").Value = "LUZ MIRYAN Y WILLIAN MOSQUERA"

# --- Enlace Técnico (M) correction ---
$ws.Range("M6").Value = "ANDREA GONZALEZ"
